# Update countries & provincias Spain
# Applies updated statistics to specific country rows in the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 29 - Singapur
$ws.Range("B29").Value = 18778
$ws.Range("C29").Value = 573
$ws.Range("E29").Value = 17352

# Row 36 - Polonia
$ws.Range("D36").Value = 4095
$ws.Range("E36").Value = 8920

# Row 37 - Rumania
$ws.Range("E37").Value = 7491
$ws.Range("G37").Value = 13
$ws.Range("H37").Value = 803

# Row 38 - Ucrania
$ws.Range("B38").Value = 12331
$ws.Range("C38").Value = 418
$ws.Range("D38").Value = 1619
$ws.Range("E38").Value = 10409
$ws.Range("F38").Value = 160
$ws.Range("G38").Value = 15
$ws.Range("H38").Value = 303

# Row 46 - Noruega
$ws.Range("E46").Value = 7603
$ws.Range("G46").Value = 1
$ws.Range("H46").Value = 212

# Row 69 - Armenia
$ws.Range("B69").Value = 2507
$ws.Range("C69").Value = 121
$ws.Range("D69").Value = 1071
$ws.Range("E69").Value = 1397
$ws.Range("G69").Value = 4
$ws.Range("H69").Value = 39

# Row 86 - Lituania
$ws.Range("B86").Value = 1419
$ws.Range("C86").Value = 9
$ws.Range("D86").Value = 638
$ws.Range("E86").Value = 735

# Row 94 - Letonia
$ws.Range("B94").Value = 896
$ws.Range("C94").Value = 17
$ws.Range("E94").Value = 532

$wb.Save()
